$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 2.05
$ws.Range("R2").Value = 1.68
$ws.Range("U2").Value = 1.77
$ws.Range("V2").Value = 1.92
$ws.Range("AT2").Value = 2.62
$ws.Range("Q3").Value = 2.15
$ws.Range("R3").Value = 1.62
$ws.Range("U3").Value = 1.92
$ws.Range("V3").Value = 1.77
$ws.Range("BD3").Value = 151
$ws.Range("J4").Value = 2.37
$ws.Range("Q4").Value = 2.25
$ws.Range("R4").Value = 1.57
$ws.Range("V4").Value = 1.58
$ws.Range("G5").Value = 2.55
$ws.Range("H5").Value = 2.75
$ws.Range("G6").Value = 1.8
$ws.Range("G7").Value = 2.4
$ws.Range("I7").Value = 2.6
$ws.Range("O7").Value = 1.22
$ws.Range("P7").Value = 4
$ws.Range("U7").Value = 1.63
$ws.Range("I8").Value = 1.83
$ws.Range("U8").Value = 1.77
$ws.Range("V8").Value = 1.87
$ws.Range("M9").Value = 1.03
$ws.Range("O9").Value = 1.22
$ws.Range("G11").Value = 1.75
$ws.Range("M11").Value = 1.02
$ws.Range("N11").Value = 17
$ws.Range("O11").Value = 1.13
$ws.Range("P11").Value = 5
$ws.Range("R11").Value = 2.35
$ws.Range("Q12").Value = 1.67
$ws.Range("U12").Value = 1.54
$ws.Range("G13").Value = 3.5
$ws.Range("H13").Value = 4.1
$ws.Range("I13").Value = 1.81
$ws.Range("L13").Value = 2.38
$ws.Range("Q13").Value = 1.4
$ws.Range("R13").Value = 2.88
$ws.Range("U13").Value = 1.37
$ws.Range("V13").Value = 2.75
$ws.Range("AC13").Value = 21
$ws.Range("AD13").Value = 8.5
$ws.Range("AE13").Value = 11
$ws.Range("AI13").Value = 13
$ws.Range("AK13").Value = 19
$ws.Range("AN13").Value = 6
$ws.Range("AO13").Value = 17
$ws.Range("AP13").Value = 19
$ws.Range("AX13").Value = 9.5
$ws.Range("I14").Value = 1.71
$ws.Range("U14").Value = 1.33
$ws.Range("G17").Value = 1.53
$ws.Range("Q17").Value = 1.88
$ws.Range("R17").Value = 1.98
$ws.Range("I18").Value = 2.88
$ws.Range("G20").Value = 1.73
$ws.Range("Q25").Value = 1.89
$ws.Range("R25").Value = 1.84
$ws.Range("U25").Value = 1.8
$ws.Range("V25").Value = 1.95
$ws.Range("U26").Value = 1.62
$ws.Range("V27").Value = 1.73
$ws.Range("U28").Value = 1.73
$ws.Range("G31").Value = 2.3
$ws.Range("I31").Value = 2.7
$ws.Range("I32").Value = 2.15
$ws.Range("G33").Value = 1.5
$ws.Range("J35").Value = 2.87
$ws.Range("Q35").Value = 1.95
$ws.Range("R35").Value = 1.85
